# Add two more plate rows ("AE" and "AF") to the 1536-well layout, extending
# the existing A2:AW32 table down to A2:AW34 (support for lowercase/extra
# row indexing beyond "AD").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32 ("AD") is the last fully populated data row; duplicate its B:AW
# well-id values (every row shares the same 1..48 well-id header strings)
# into the two new rows.
$ws.Range("B32:AW32").Copy() | Out-Null
$ws.Range("B33:AW33").PasteSpecial(-4104) | Out-Null   # xlPasteAll
$ws.Range("B32:AW32").Copy() | Out-Null
$ws.Range("B34:AW34").PasteSpecial(-4104) | Out-Null   # xlPasteAll

# Fill in the new row-letter labels in column A.
$ws.Range("A33").Value2 = "AE"
$ws.Range("A34").Value2 = "AF"

# Apply the same header formatting (white text on gray fill) used by the
# rest of column A onto the two new label cells.
$ws.Range("A32").Copy() | Out-Null
$ws.Range("A33:A34").PasteSpecial(-4122) | Out-Null    # xlPasteFormats

$excel.CutCopyMode = 0

# Move the active selection down past the newly added rows.
$ws.Range("A35").Select() | Out-Null
